$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data after the edit: the first 3 matches (RCD Mallorca, Villarreal CF,
# FC Sheriff Tiraspol) were removed, and 3 new matches were appended at
# the end of the fixture list (Cadiz CF, Valencia CF, Elche CF).
$data = @(
    @("Real Madrid CF - Athletic Club de Bilbao", "17/10/2021"),
    @("Real Madrid CF - Osasuna", "27/10/2021"),
    @("Real Madrid CF - FC Shakhtar Donetsk", "03/11/2021"),
    @("Real Madrid CF - Rayo Vallecano", "07/11/2021"),
    @("Real Madrid CF - Sevilla FC", "28/11/2021"),
    @("Real Madrid CF - Inter Milan", "07/12/2021"),
    @("Real Madrid CF - Atlético de Madrid", "12/12/2021"),
    @("Real Madrid CF - Cadiz CF", "19/12/2021"),
    @("Real Madrid CF - Valencia CF", "09/01/2022"),
    @("Real Madrid CF - Elche CF", "23/01/2022")
)

# Temporarily force column B to Text format so the date-like strings
# (e.g. "03/11/2021") are not auto-parsed into date serial numbers when
# assigned below - the original workbook stores these as plain text
# shared strings, not as dates.
$ws.Range("B1:B10").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value2 = $data[$i][0]
    $ws.Cells.Item($row, 2).Value2 = $data[$i][1]
}

# Restore the cells to the workbook's original (default/General) style so
# no stray per-cell formatting is left behind.
$ws.Range("B1:B10").ClearFormats()
